$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3798.9812
$ws.Range("I64").Value = 3602.3704
$ws.Range("J64").Value = 4003.1538
$ws.Range("K64").Value = 3602.3704
$ws.Range("L64").Value = 4003.1538
$ws.Range("M64").Value = -3354.3704
$ws.Range("N64").Value = -4499.1538
$ws.Range("H67").Value = 3798.9812
$ws.Range("I67").Value = 3602.3704
$ws.Range("J67").Value = 4003.1538
$ws.Range("K67").Value = 3602.3704
$ws.Range("L67").Value = 4003.1538
$ws.Range("M67").Value = -2744.3704
$ws.Range("N67").Value = -5719.1538
$ws.Range("H68").Value = 37633.332
$ws.Range("J68").Value = 37633.332
$ws.Range("L68").Value = 37633.332
$ws.Range("N68").Value = -39131.332
$ws.Range("H71").Value = 37633.332
$ws.Range("J71").Value = 37633.332
$ws.Range("L71").Value = 112899.996
$ws.Range("N71").Value = -120387.996
$ws.Range("H113").Value = 1865.2593
$ws.Range("I113").Value = 1767.4286
$ws.Range("J113").Value = 1899.5
$ws.Range("K113").Value = 1767.4286
$ws.Range("L113").Value = 1899.5
$ws.Range("M113").Value = 1486.5714
$ws.Range("N113").Value = -8407.5
$ws.Range("H132").Value = 1986.5
$ws.Range("I132").Value = 2135.0667
$ws.Range("K132").Value = 6405.2001
$ws.Range("M132").Value = -3875.2001
$ws.Range("H138").Value = 4185.5845
$ws.Range("I138").Value = 1217.2727
$ws.Range("J138").Value = 4790.2407
$ws.Range("K138").Value = 3651.8181
$ws.Range("L138").Value = 14370.7221
$ws.Range("M138").Value = 1488.1819
$ws.Range("N138").Value = -24650.7221

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7506.143
$ws.Range("I2").Value = 7506.143
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 7506.143
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -7393.143
$ws.Range("N2").ClearContents()
$ws.Range("H32").Value = 27031764
$ws.Range("I32").Value = 29414846
$ws.Range("J32").Value = 23504.666
$ws.Range("K32").Value = 29414846
$ws.Range("L32").Value = 23504.666
$ws.Range("M32").Value = -29414559
$ws.Range("N32").Value = -24078.666
$ws.Range("H102").Value = 1877
$ws.Range("I102").Value = 2002.6666
$ws.Range("K102").Value = 2002.6666
$ws.Range("M102").Value = -380.6666
$ws.Range("H116").Value = 7506.143
$ws.Range("I116").Value = 7506.143
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 7506.143
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -5212.143
$ws.Range("N116").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7506.143
$ws.Range("I3").Value = 7506.143
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 7506.143
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -7392.143
$ws.Range("N3").ClearContents()
$ws.Range("H86").Value = 1906.0588
$ws.Range("I86").Value = 1839.7333
$ws.Range("J86").Value = 2403.5
$ws.Range("K86").Value = 1839.7333
$ws.Range("L86").Value = 2403.5
$ws.Range("M86").Value = -716.7333000000001
$ws.Range("N86").Value = -4649.5
$ws.Range("H89").Value = 1906.0588
$ws.Range("I89").Value = 1839.7333
$ws.Range("J89").Value = 2403.5
$ws.Range("K89").Value = 9198.666500000001
$ws.Range("L89").Value = 12017.5
$ws.Range("M89").Value = -3582.666500000001
$ws.Range("N89").Value = -23249.5
$ws.Range("H94").Value = 887.6429000000001
$ws.Range("I94").Value = 893.9167
$ws.Range("J94").Value = 850
$ws.Range("K94").Value = 893.9167
$ws.Range("L94").Value = 850
$ws.Range("M94").Value = -442.9167
$ws.Range("N94").Value = -1752
$ws.Range("H99").Value = 1711.8846
$ws.Range("I99").Value = 1267.65
$ws.Range("J99").Value = 3192.6667
$ws.Range("K99").Value = 1267.65
$ws.Range("L99").Value = 3192.6667
$ws.Range("M99").Value = 230.3499999999999
$ws.Range("N99").Value = -6188.6667
$ws.Range("H105").Value = 2781.1633
$ws.Range("I105").Value = 1645.5555
$ws.Range("J105").Value = 2896
$ws.Range("K105").Value = 1645.5555
$ws.Range("L105").Value = 2896
$ws.Range("M105").Value = 101.4445000000001
$ws.Range("N105").Value = -6390

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2574
$ws.Range("I62").Value = 2512.111
$ws.Range("J62").Value = 2666.8333
$ws.Range("K62").Value = 2512.111
$ws.Range("L62").Value = 2666.8333
$ws.Range("M62").Value = -1888.111
$ws.Range("N62").Value = -3914.8333
$ws.Range("H65").Value = 2574
$ws.Range("I65").Value = 2512.111
$ws.Range("J65").Value = 2666.8333
$ws.Range("K65").Value = 12560.555
$ws.Range("L65").Value = 13334.1665
$ws.Range("M65").Value = -9440.555
$ws.Range("N65").Value = -19574.1665
$ws.Range("H99").Value = 2346.8333
$ws.Range("I99").Value = 1965.5714
$ws.Range("J99").Value = 2880.6
$ws.Range("K99").Value = 1965.5714
$ws.Range("L99").Value = 2880.6
$ws.Range("M99").Value = -467.5714
$ws.Range("N99").Value = -5876.6
$ws.Range("H105").Value = 1779
$ws.Range("I105").Value = 1662
$ws.Range("J105").Value = 1974
$ws.Range("K105").Value = 1662
$ws.Range("L105").Value = 1974
$ws.Range("M105").Value = 85
$ws.Range("N105").Value = -5468
$ws.Range("H126").Value = 2346.8333
$ws.Range("I126").Value = 1965.5714
$ws.Range("J126").Value = 2880.6
$ws.Range("K126").Value = 5896.7142
$ws.Range("L126").Value = 8641.799999999999
$ws.Range("M126").Value = -3426.7142
$ws.Range("N126").Value = -13581.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 848.15625
$ws.Range("I131").Value = 519
$ws.Range("J131").Value = 909.1111
$ws.Range("K131").Value = 1557
$ws.Range("L131").Value = 2727.3333
$ws.Range("M131").Value = 3483
$ws.Range("N131").Value = -12807.3333

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2825.8057
$ws.Range("I80").Value = 2752.037
$ws.Range("J80").Value = 3047.111
$ws.Range("K80").Value = 2752.037
$ws.Range("L80").Value = 3047.111
$ws.Range("M80").Value = -1754.037
$ws.Range("N80").Value = -5043.111
$ws.Range("H83").Value = 2825.8057
$ws.Range("I83").Value = 2752.037
$ws.Range("J83").Value = 3047.111
$ws.Range("K83").Value = 13760.185
$ws.Range("L83").Value = 15235.555
$ws.Range("M83").Value = -8768.184999999999
$ws.Range("N83").Value = -25219.555

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 536.5
$ws.Range("I22").Value = 506.55554
$ws.Range("J22").Value = 590.4
$ws.Range("K22").Value = 506.55554
$ws.Range("L22").Value = 590.4
$ws.Range("M22").Value = -211.55554
$ws.Range("N22").Value = -1180.4
$ws.Range("H27").Value = 536.5
$ws.Range("I27").Value = 506.55554
$ws.Range("J27").Value = 590.4
$ws.Range("K27").Value = 506.55554
$ws.Range("L27").Value = 590.4
$ws.Range("M27").Value = -399.55554
$ws.Range("N27").Value = -804.4
$ws.Range("H40").Value = 3860.923
$ws.Range("I40").Value = 3566.45
$ws.Range("J40").Value = 4842.5
$ws.Range("K40").Value = 3566.45
$ws.Range("L40").Value = 4842.5
$ws.Range("M40").Value = -3430.45
$ws.Range("N40").Value = -5114.5
$ws.Range("H68").Value = 1814.7273
$ws.Range("I68").Value = 1666.8889
$ws.Range("J68").Value = 2480
$ws.Range("K68").Value = 1666.8889
$ws.Range("L68").Value = 2480
$ws.Range("M68").Value = -917.8888999999999
$ws.Range("N68").Value = -3978
$ws.Range("H71").Value = 1814.7273
$ws.Range("I71").Value = 1666.8889
$ws.Range("J71").Value = 2480
$ws.Range("K71").Value = 8334.4445
$ws.Range("L71").Value = 12400
$ws.Range("M71").Value = -4590.4445
$ws.Range("N71").Value = -19888
$ws.Range("H82").Value = 2266.6667
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 2400
$ws.Range("K82").Value = 2000
$ws.Range("L82").Value = 2400
$ws.Range("M82").Value = -1639
$ws.Range("N82").Value = -3122
$ws.Range("H85").Value = 2266.6667
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 2400
$ws.Range("K85").Value = 2000
$ws.Range("L85").Value = 2400
$ws.Range("M85").Value = -752
$ws.Range("N85").Value = -4896
$ws.Range("H93").Value = 8559.6
$ws.Range("I93").Value = 10249.4375
$ws.Range("J93").Value = 1800.25
$ws.Range("K93").Value = 10249.4375
$ws.Range("L93").Value = 1800.25
$ws.Range("M93").Value = -9001.4375
$ws.Range("N93").Value = -4296.25
$ws.Range("H100").Value = 39232310
$ws.Range("I100").Value = 910778.0600000001
$ws.Range("K100").Value = 910778.0600000001
$ws.Range("M100").Value = -910237.0600000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11124133
$ws.Range("I62").Value = 16684700
$ws.Range("K62").Value = 16684700
$ws.Range("M62").Value = -16684076
$ws.Range("H63").Value = 25000
$ws.Range("J63").Value = 25000
$ws.Range("L63").Value = 25000
$ws.Range("N63").Value = -26248
$ws.Range("H65").Value = 11124133
$ws.Range("I65").Value = 16684700
$ws.Range("K65").Value = 83423500
$ws.Range("M65").Value = -83420380
$ws.Range("H66").Value = 25000
$ws.Range("J66").Value = 25000
$ws.Range("L66").Value = 75000
$ws.Range("N66").Value = -81240
$ws.Range("H69").Value = 16180
$ws.Range("J69").Value = 16180
$ws.Range("L69").Value = 16180
$ws.Range("N69").Value = -17678
$ws.Range("H72").Value = 16180
$ws.Range("J72").Value = 16180
$ws.Range("L72").Value = 48540
$ws.Range("N72").Value = -56028
$ws.Range("H81").Value = 1943.75
$ws.Range("J81").Value = 1943.75
$ws.Range("L81").Value = 3887.5
$ws.Range("N81").Value = -6009.5
$ws.Range("H84").Value = 1943.75
$ws.Range("J84").Value = 1943.75
$ws.Range("L84").Value = 19437.5
$ws.Range("N84").Value = -30045.5
$ws.Range("H107").Value = 363.72726
$ws.Range("I107").Value = 363.72726
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1091.18178
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 828.8182200000001
$ws.Range("N107").ClearContents()
